$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1752.2755
$ws.Range("I15").Value = 1752.2755
$ws.Range("K15").Value = 5256.8265
$ws.Range("M15").Value = -5087.8265

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1557.9375
$ws.Range("I40").Value = 1362.5
$ws.Range("J40").Value = 2144.25
$ws.Range("K40").Value = 1362.5
$ws.Range("L40").Value = 2144.25
$ws.Range("M40").Value = -1187.5
$ws.Range("N40").Value = -2494.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3689.2144
$ws.Range("I64").Value = 3039.8
$ws.Range("J64").Value = 4050
$ws.Range("K64").Value = 3039.8
$ws.Range("L64").Value = 4050
$ws.Range("M64").Value = -2791.8
$ws.Range("N64").Value = -4546

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3689.2144
$ws.Range("I67").Value = 3039.8
$ws.Range("J67").Value = 4050
$ws.Range("K67").Value = 3039.8
$ws.Range("L67").Value = 4050
$ws.Range("M67").Value = -2181.8
$ws.Range("N67").Value = -5766

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 7918.3335
$ws.Range("I86").Value = 1393.3
$ws.Range("J86").Value = 20968.4
$ws.Range("K86").Value = 1393.3
$ws.Range("L86").Value = 20968.4
$ws.Range("M86").Value = -270.3
$ws.Range("N86").Value = -23214.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1500
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1500
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1500
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -2312

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 7918.3335
$ws.Range("I89").Value = 1393.3
$ws.Range("J89").Value = 20968.4
$ws.Range("K89").Value = 6966.5
$ws.Range("L89").Value = 104842
$ws.Range("M89").Value = -1350.5
$ws.Range("N89").Value = -116074

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1500
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1500
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1500
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -4308

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 38235.4
$ws.Range("J105").Value = 38235.4
$ws.Range("L105").Value = 38235.4
$ws.Range("N105").Value = -45223.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 845.3182
$ws.Range("I107").Value = 568.4737
$ws.Range("J107").Value = 2598.6667
$ws.Range("K107").Value = 568.4737
$ws.Range("L107").Value = 2598.6667
$ws.Range("M107").Value = 1351.5263
$ws.Range("N107").Value = -6438.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 420.44446
$ws.Range("I115").Value = 420.44446
$ws.Range("K115").Value = 1261.33338
$ws.Range("M115").Value = 305.66662

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4709.6665
$ws.Range("J116").Value = 5964.6
$ws.Range("L116").Value = 5964.6
$ws.Range("N116").Value = -12848.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 244862.92
$ws.Range("J129").Value = 313668.8
$ws.Range("L129").Value = 941006.3999999999
$ws.Range("N129").Value = -951006.3999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2499.2222
$ws.Range("I132").Value = 2545.4187
$ws.Range("J132").Value = 1506
$ws.Range("K132").Value = 7636.256100000001
$ws.Range("L132").Value = 4518
$ws.Range("M132").Value = -5106.256100000001
$ws.Range("N132").Value = -9578

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3298.7778
$ws.Range("I141").Value = 2922.25
$ws.Range("J141").Value = 3600
$ws.Range("K141").Value = 8766.75
$ws.Range("L141").Value = 10800
$ws.Range("M141").Value = -3586.75
$ws.Range("N141").Value = -21160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10398.904
$ws.Range("I32").Value = 7564.8486
$ws.Range("J32").Value = 20790.445
$ws.Range("K32").Value = 7564.8486
$ws.Range("L32").Value = 20790.445
$ws.Range("M32").Value = -7277.8486
$ws.Range("N32").Value = -21364.445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 25642312
$ws.Range("I74").Value = 34483124
$ws.Range("K74").Value = 34483124
$ws.Range("M74").Value = -34482250

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 25642312
$ws.Range("I77").Value = 34483124
$ws.Range("K77").Value = 172415620
$ws.Range("M77").Value = -172411252

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 41667380
$ws.Range("I97").Value = 412.9
$ws.Range("K97").Value = 412.9
$ws.Range("M97").Value = 83.10000000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 721.46155
$ws.Range("I110").Value = 620.2222
$ws.Range("J110").Value = 949.25
$ws.Range("K110").Value = 620.2222
$ws.Range("L110").Value = 949.25
$ws.Range("M110").Value = 1424.7778
$ws.Range("N110").Value = -5039.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 54957.25
$ws.Range("J140").Value = 54957.25
$ws.Range("L140").Value = 54957.25
$ws.Range("N140").Value = -65317.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 29354
$ws.Range("J87").Value = 29354
$ws.Range("L87").Value = 29354
$ws.Range("N87").Value = -31850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H90").Value = 29354
$ws.Range("J90").Value = 29354
$ws.Range("L90").Value = 88062
$ws.Range("N90").Value = -100542

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 822.6923
$ws.Range("J107").Value = 828.25
$ws.Range("L107").Value = 828.25
$ws.Range("N107").Value = -4668.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 19426.072
$ws.Range("I58").Value = 1486.5333
$ws.Range("J58").Value = 40125.54
$ws.Range("K58").Value = 1486.5333
$ws.Range("L58").Value = 40125.54
$ws.Range("M58").Value = -1283.5333
$ws.Range("N58").Value = -40531.54

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 19426.072
$ws.Range("I136").Value = 1486.5333
$ws.Range("J136").Value = 40125.54
$ws.Range("K136").Value = 4459.5999
$ws.Range("L136").Value = 120376.62
$ws.Range("M136").Value = -1909.5999
$ws.Range("N136").Value = -125476.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 19.866667
$ws.Range("J2").Value = 9
$ws.Range("L2").Value = 54
$ws.Range("N2").Value = -280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1743.3
$ws.Range("J5").Value = 2611
$ws.Range("L5").Value = 7833
$ws.Range("N5").Value = -8057

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 843.2857
$ws.Range("I34").Value = 399.5
$ws.Range("J34").Value = 1020.8
$ws.Range("K34").Value = 1198.5
$ws.Range("L34").Value = 3062.4
$ws.Range("M34").Value = -1114.5
$ws.Range("N34").Value = -3230.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 145
$ws.Range("I60").Value = 145
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 435
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -184
$ws.Range("N60").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1104.5454
$ws.Range("J121").Value = 1142.8572
$ws.Range("L121").Value = 3428.5716
$ws.Range("N121").Value = -6048.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 114416.555
$ws.Range("I131").Value = 665
$ws.Range("J131").Value = 117061.945
$ws.Range("K131").Value = 1995
$ws.Range("L131").Value = 351185.835
$ws.Range("M131").Value = 3045
$ws.Range("N131").Value = -361265.835

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1743.3
$ws.Range("J135").Value = 2611
$ws.Range("L135").Value = 23499
$ws.Range("N135").Value = -28569

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4120.95
$ws.Range("J80").Value = 4584.5386
$ws.Range("L80").Value = 4584.5386
$ws.Range("N80").Value = -6580.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4120.95
$ws.Range("J83").Value = 4584.5386
$ws.Range("L83").Value = 22922.693
$ws.Range("N83").Value = -32906.693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3077181.8
$ws.Range("I107").Value = 269.875
$ws.Range("J107").Value = 8547247
$ws.Range("K107").Value = 269.875
$ws.Range("L107").Value = 8547247
$ws.Range("M107").Value = 1650.125
$ws.Range("N107").Value = -8551087

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 49952.57
$ws.Range("J135").Value = 49952.57
$ws.Range("L135").Value = 49952.57
$ws.Range("N135").Value = -60092.57

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2011
$ws.Range("I22").Value = 2780
$ws.Range("J22").Value = 582.8570999999999
$ws.Range("K22").Value = 2780
$ws.Range("L22").Value = 582.8570999999999
$ws.Range("M22").Value = -2485
$ws.Range("N22").Value = -1172.8571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2011
$ws.Range("I27").Value = 2780
$ws.Range("J27").Value = 582.8570999999999
$ws.Range("K27").Value = 2780
$ws.Range("L27").Value = 582.8570999999999
$ws.Range("M27").Value = -2673
$ws.Range("N27").Value = -796.8570999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3765.2856
$ws.Range("I40").Value = 2894.4707
$ws.Range("J40").Value = 7466.25
$ws.Range("K40").Value = 2894.4707
$ws.Range("L40").Value = 7466.25
$ws.Range("M40").Value = -2758.4707
$ws.Range("N40").Value = -7738.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4485.3335
$ws.Range("I61").Value = 1882.4
$ws.Range("K61").Value = 1882.4
$ws.Range("M61").Value = -1680.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2533.8
$ws.Range("I68").Value = 2600
$ws.Range("J68").Value = 2505.4285
$ws.Range("K68").Value = 2600
$ws.Range("L68").Value = 2505.4285
$ws.Range("M68").Value = -1851
$ws.Range("N68").Value = -4003.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2533.8
$ws.Range("I71").Value = 2600
$ws.Range("J71").Value = 2505.4285
$ws.Range("K71").Value = 13000
$ws.Range("L71").Value = 12527.1425
$ws.Range("M71").Value = -9256
$ws.Range("N71").Value = -20015.1425

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1205.1111
$ws.Range("I93").Value = 1042.75
$ws.Range("J93").Value = 2504
$ws.Range("K93").Value = 1042.75
$ws.Range("L93").Value = 2504
$ws.Range("M93").Value = 205.25
$ws.Range("N93").Value = -5000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4485.3335
$ws.Range("I113").Value = 1882.4
$ws.Range("K113").Value = 1882.4
$ws.Range("M113").Value = 287.5999999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 504198.66
$ws.Range("I132").Value = 1005115.94
$ws.Range("K132").Value = 3015347.82
$ws.Range("M132").Value = -3012817.82

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 30000000
$ws.Range("J70").Value = 30000000
$ws.Range("L70").Value = 30000000
$ws.Range("N70").Value = -30000630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 30000000
$ws.Range("J73").Value = 30000000
$ws.Range("L73").Value = 30000000
$ws.Range("N73").Value = -30002184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 55024468
$ws.Range("I107").Value = 90909260
$ws.Range("K107").Value = 272727780
$ws.Range("M107").Value = -272725860
